# B3-and-B4-PowerPoint.pptx edit
#
# The canonical OOXML diff swaps the contents of ppt/theme/theme1.xml
# (the slide-master theme, originally the "Integral" / "Red Violet"
# palette) and ppt/theme/theme2.xml (the notes-master theme, originally
# the default "Office Theme" palette) - i.e. after the edit the slide
# master uses the Office Theme colours and the notes master uses the
# Integral/Red Violet colours. Font scheme and format scheme are
# identical between the two themes, so the only real content change is
# the 12 theme colours (and the cosmetic theme/colour-scheme names).
#
# Apply it the way PowerPoint exposes theme colours: through the
# SlideMaster's Theme.ThemeColorScheme collection (ColorFormat.RGB is
# the supported read/write surface for theme colours in the object
# model).

function Convert-HexToComRgb([string]$hex) {
    $val = [Convert]::ToInt32($hex, 16)
    $r = ($val -shr 16) -band 0xFF
    $g = ($val -shr 8) -band 0xFF
    $b = $val -band 0xFF
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target palette for the slide-master theme (the Office Theme colours
# that theme2.xml held before the edit, and that theme1.xml holds
# after it). Order matches ThemeColorScheme's fixed 12-slot layout:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = Convert-HexToComRgb $officeThemeColors[$i - 1]
}
